$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $c = $sheet.Range($addr)
    $c.Formula = "'" + $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '26.709.49'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '1.632.95'
$ws.Range('E3').Value = '  +1.50%  '
$ws.Range('E4').Value = '  -0.06%  '
Set-TextValue $ws 'D5' '213.59'
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('E7').Value = '  +0.93%  '
$ws.Range('E8').Value = '  +0.61%  '
$ws.Range('E9').Value = '  +0.88%  '
Set-TextValue $ws 'D10' '19.06'
Set-TextValue $ws 'D11' '0.0833'
$ws.Range('E11').Value = '  +2.32%  '
$ws.Range('D12').Value = '1.857.23'
$ws.Range('E12').Value = '  +1.33%  '
$ws.Range('D13').Value = '1.626.94'
$ws.Range('E13').Value = '  +1.03%  '
$ws.Range('E14').Value = '  +0.23%  '
Set-TextValue $ws 'D15' '0.525'
$ws.Range('E15').Value = '  +2.01%  '
$ws.Range('D16').Value = '26.670.72'
$ws.Range('E16').Value = '  +1.47%  '
Set-TextValue $ws 'D17' '63.07'
$ws.Range('E17').Value = '  +2.40%  '
$ws.Range('D18').Value = '0.0₃0734'
$ws.Range('E18').Value = '  +0.67%  '
Set-TextValue $ws 'D19' '209.41'
$ws.Range('E19').Value = '  +2.93%  '
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('E21').Value = '  +1.03%  '
Set-TextValue $ws 'D22' '9.40'
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('E23').Value = '  +1.60%  '
$ws.Range('E24').Value = '  -1.69%  '
Set-TextValue $ws 'D25' '145.75'
$ws.Range('E25').Value = '  +0.89%  '
$ws.Range('E26').Value = '  -0.13%  '
Set-TextValue $ws 'D27' '0.121'
$ws.Range('E27').Value = '  -1.21%  '
Set-TextValue $ws 'D28' '15.36'
$ws.Range('E28').Value = '  +1.00%  '
Set-TextValue $ws 'D29' '6.66'
$ws.Range('E29').Value = '  +1.27%  '
Set-TextValue $ws 'D30' '0.0519'
$ws.Range('E30').Value = '  +6.34%  '
$ws.Range('E31').Value = '  +0.34%  '
$ws.Range('E32').Value = '  +1.30%  '
$ws.Range('E33').Value = '  +0.93%  '
$ws.Range('E34').Value = '  +1.58%  '
$ws.Range('E35').Value = '  -0.37%  '
$ws.Range('D36').Value = '1.165.80'
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('E37').Value = '  +0.72%  '
Set-TextValue $ws 'D38' '0.814'
$ws.Range('E38').Value = '  +2.22%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  -0.63%  '
Set-TextValue $ws 'D41' '0.503'
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('E42').Value = '  +3.22%  '
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('D44').Value = '1.767.71'
$ws.Range('E44').Value = '  +1.32%  '
Set-TextValue $ws 'D45' '92.35'
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('E46').Value = '  +1.88%  '
Set-TextValue $ws 'D47' '54.68'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('E48').Value = '  +1.18%  '
Set-TextValue $ws 'D49' '7.72'
$ws.Range('E49').Value = '  +6.79%  '
Set-TextValue $ws 'D50' '0.410'
$ws.Range('E50').Value = '  +0.83%  '
$ws.Range('E51').Value = '  -0.03%  '

Write-Host "Updated cryptos list"